# Auto-generated edit script: update crypto price/volume table cells
# to match the latest scraped data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.456.10'
$ws.Range('E2').Value = '  +3.81%  '
$ws.Range('D3').Value = '3.618.06'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '628.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.61'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.13%  '
$ws.Range('D7').Value = '3.616.76'
$ws.Range('E7').Value = '  +2.78%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.63%  '
$ws.Range('E10').Value = '  +4.98%  '
$ws.Range('E11').Value = '  +6.62%  '
$ws.Range('E12').Value = '  +3.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000225'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.63%  '
$ws.Range('D15').Value = '4.233.00'
$ws.Range('E15').Value = '  +2.87%  '
$ws.Range('D16').Value = '3.618.38'
$ws.Range('E16').Value = '  +2.73%  '
$ws.Range('D17').Value = '69.065.28'
$ws.Range('E17').Value = '  +3.21%  '
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +10.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '459.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.640'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.47%  '
$ws.Range('E25').Value = '  +13.52%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.88%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.763.52'
$ws.Range('E27').Value = '  +2.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +13.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.62'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.99%  '
$ws.Range('E32').Value = '  +12.04%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.60'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.29%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.995'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.53%  '
$ws.Range('E35').Value = '  +6.73%  '
$ws.Range('E36').Value = '  +3.61%  '
$ws.Range('D37').Value = '3.608.87'
$ws.Range('E37').Value = '  +2.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.41'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.72%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0928'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.44%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '176.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.62'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '32.07'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +17.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.911'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.58%  '
$ws.Range('E47').Value = '  +13.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.81'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.89%  '
$ws.Range('E51').Value = '  +8.01%  '
